$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("E2").Value = "2024.04.11 10:00 - 04.12 17:00"
$ws1.Range("F2").Value = 840
$ws1.Range("E3").Value = "2024.05.01 09:30 - 05.02 17:30"
$ws1.Range("F3").Value = 4313
$ws1.Range("E4").Value = "2024.05.19 10:00 - 05.19 17:00"
$ws1.Range("F4").Value = 121
$ws1.Range("E5").Value = "2024.06.09 10:00 - 06.10 17:00"
$ws1.Range("F5").Value = 772

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("E2").Value = "2024.06.22 20:00 - 06.22 21:30"

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("E2").Value = "2024.04.11 10:00 - 04.12 17:00"
$ws4.Range("F2").Value = 840
$ws4.Range("E3").Value = "2024.05.01 09:30 - 05.02 17:30"
$ws4.Range("F3").Value = 4313
$ws4.Range("E4").Value = "2024.05.19 10:00 - 05.19 17:00"
$ws4.Range("F4").Value = 121
$ws4.Range("E5").Value = "2024.06.09 10:00 - 06.10 17:00"
$ws4.Range("F5").Value = 772
$ws4.Range("E6").Value = "2024.06.22 20:00 - 06.22 21:30"
